$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lowercase the header labels
$ws.Range("A1").Value = "matricula"
$ws.Range("B1").Value = "valor"

# Apply the same cell formatting already used by column A's data rows
# (A2:A4) to the header row and to the "valor" data column, so the whole
# table shares one consistent style.
$headerStyle = $ws.Range("A2").Style
$ws.Range("A1").Style = $headerStyle
$ws.Range("B1").Style = $headerStyle
$ws.Range("B2").Style = $headerStyle
$ws.Range("B3").Style = $headerStyle
$ws.Range("B4").Style = $headerStyle

# Move the active selection to B2
$ws.Range("B2").Select() | Out-Null
